$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JatMagus")

$ws.Range("A32").Value = "etzui"
$ws.Range("B32").Value = "Mágus"
$ws.Range("C32").Value = "föld botja"
$ws.Range("D32").Value = "/Images/Karakterek/magus1.png"

$ws.Range("A31:D31").Copy()
$ws.Range("A32:D32").PasteSpecial(-4122)
